$d = $word.ActiveDocument

# The target paragraph currently holds the markup "<id>p072v_1</id>" split
# across three runs: "<id>" (Courier New, gold), "p072v_1" (plain black),
# "</id>" (Courier New, gold). We collapse it into a single run -
# "<id>p072v_1</id>" - using the formatting of the first ("<id>") run.

$paraRange = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<id>p072v_1</id>*") {
        $paraRange = $p.Range.Duplicate
        break
    }
}
if ($paraRange -eq $null) {
    throw "Could not locate paragraph containing '<id>p072v_1</id>'"
}

# Find the "<id>" run inside the paragraph; keep its Range (and formatting).
$tagRange = $paraRange.Duplicate
$found = $tagRange.Find.Execute("<id>", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate '<id>' run inside target paragraph"
}

# Everything after "<id>" up to (not including) the paragraph mark is the
# remainder of the markup ("p072v_1</id>") spread across the other two runs.
$tailRange = $d.Range($tagRange.End, $paraRange.End - 1)

# Delete those trailing runs, then append their combined text onto the
# "<id>" run so the whole string lives in one run with that run's formatting.
$tailText = $tailRange.Text
$tailRange.Delete()
$tagRange.InsertAfter($tailText)
